$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the previous year row (A4) onto the new label cell A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 11.8
$ws.Range("C5").Value = 18.1
$ws.Range("D5").Value = -88.40000000000001
$ws.Range("E5").Value = 85
$ws.Range("F5").Value = 6.5
$ws.Range("G5").Value = "'"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 44.8
$ws.Range("I5").Value = -12.6
$ws.Range("J5").Value = -10.9
$ws.Range("K5").Value = -56
$ws.Range("L5").Value = -49.1
$ws.Range("M5").Value = 55.3
$ws.Range("N5").Value = -28.1
$ws.Range("O5").Value = -18.4
$ws.Range("P5").Value = -48.5
$ws.Range("Q5").Value = -26.8
$ws.Range("R5").Value = "'"
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").Value = -50.5
$ws.Range("T5").Value = 0.6
$ws.Range("U5").Value = 21.5
$ws.Range("V5").Value = -16.5
$ws.Range("W5").Value = 6.4
$ws.Range("X5").Value = 6.7
$ws.Range("Y5").Value = 23.8
$ws.Range("Z5").Value = 5.6
$ws.Range("AA5").Value = 21.4
$ws.Range("AB5").Value = -8
$ws.Range("AC5").Value = -8.800000000000001
$ws.Range("AD5").Value = -0.4
$ws.Range("AE5").Value = 12.3
$ws.Range("AF5").Value = -16.3
$ws.Range("AG5").Value = "'"
$ws.Range("AG5").Style = "Normal"
$ws.Range("AH5").Value = -49.7
$ws.Range("AI5").Value = 68.3
$ws.Range("AJ5").Value = -88.5
$ws.Range("AK5").Value = 104.2
$ws.Range("AL5").Value = -21
$ws.Range("AM5").Value = -26.1
$ws.Range("AN5").Value = -24.9
$ws.Range("AO5").Value = -13.9
$ws.Range("AP5").Value = 11
$ws.Range("AQ5").Value = -3.2
$ws.Range("AR5").Value = 26.1
$ws.Range("AS5").Value = "'"
$ws.Range("AS5").Style = "Normal"
$ws.Range("AT5").Value = "'"
$ws.Range("AT5").Style = "Normal"
$ws.Range("AU5").Value = 1.6
$ws.Range("AV5").Value = -11.7
$ws.Range("AW5").Value = -86.3
$ws.Range("AX5").Value = -26.5
$ws.Range("AY5").Value = -6.7
$ws.Range("AZ5").Value = 7.2
$ws.Range("BA5").Value = -21.3
$ws.Range("BB5").Value = 4.1
$ws.Range("BC5").Value = "'"
$ws.Range("BC5").Style = "Normal"
$ws.Range("BD5").Value = "'"
$ws.Range("BD5").Style = "Normal"
$ws.Range("BE5").Value = -23.7
$ws.Range("BF5").Value = 13.6
$ws.Range("BG5").Value = 16.8
$ws.Range("BH5").Value = -95.40000000000001
$ws.Range("BI5").Value = 151.4
$ws.Range("BJ5").Value = 10.3
$ws.Range("BK5").Value = 72.2
$ws.Range("BL5").Value = -14.9
$ws.Range("BM5").Value = 9
$ws.Range("BN5").Value = -10
$ws.Range("BO5").Value = -25.4
$ws.Range("BP5").Value = 200.2
$ws.Range("BQ5").Value = "'"
$ws.Range("BQ5").Style = "Normal"
$ws.Range("BR5").Value = -32.1
$ws.Range("BS5").Value = 34
$ws.Range("BT5").Value = -16
$ws.Range("BU5").Value = 54.4
$ws.Range("BV5").Value = 12.2
$ws.Range("BW5").Value = 13.9
$ws.Range("BX5").Value = 44.2
$ws.Range("BY5").Value = 22.9
$ws.Range("BZ5").Value = 12.5
$ws.Range("CA5").Value = 3.8
$ws.Range("CB5").Value = 2
$ws.Range("CC5").Value = -12.7
$ws.Range("CD5").Value = "'"
$ws.Range("CD5").Style = "Normal"
$ws.Range("CE5").Value = -3.7
$ws.Range("CF5").Value = 12.8
$ws.Range("CG5").Value = 43.2
$ws.Range("CH5").Value = 225.7
$ws.Range("CI5").Value = 22
$ws.Range("CJ5").Value = -38.5
$ws.Range("CK5").Value = 7.5
$ws.Range("CL5").Value = 17.9
$ws.Range("CM5").Value = -68.5
$ws.Range("CN5").Value = 54.4
$ws.Range("CO5").Value = -9.300000000000001
$ws.Range("CP5").Value = -12.9
$ws.Range("CQ5").Value = -21.4
$ws.Range("CR5").Value = -55.1
$ws.Range("CS5").Value = 36.1
$ws.Range("CT5").Value = -2.6
$ws.Range("CU5").Value = 56.2
$ws.Range("CV5").Value = -4.5
$ws.Range("CW5").Value = -85.2
$ws.Range("CX5").Value = -7.6
$ws.Range("CY5").Value = 8.6
$ws.Range("CZ5").Value = 172.8
$ws.Range("DA5").Value = -17.5
$ws.Range("DB5").Value = -42
$ws.Range("DC5").Value = 26.7
$ws.Range("DD5").Value = 47.1
$ws.Range("DE5").Value = 20.2
$ws.Range("DF5").Value = 21.6
$ws.Range("DG5").Value = 164.9
$ws.Range("DH5").Value = 23.5
$ws.Range("DI5").Value = 24
$ws.Range("DJ5").Value = -17.8
$ws.Range("DK5").Value = 153.8
